# Update cryptocurrency price/volume data to reflect the latest market snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.835.69"
$ws.Range("E2").Value = "  +1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.090.12"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.20"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.11"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.085.33"
$ws.Range("E8").Value = "  +0.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.53"
$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("E13").Value = "  +4.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.77"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.588.96"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.926.00"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.113"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.088.62"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.65"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.18"
$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -0.59%  "

$ws.Range("E23").Value = "  -1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.88"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.38"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -3.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.33"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.91"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("E32").Value = "  +2.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.05"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("E34").Value = "  -5.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "493.46"
$ws.Range("E35").Value = "  -4.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").Value = "  +5.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.01"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.253.06"
$ws.Range("E38").Value = "  +5.65%  "

$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0800"
$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.10"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.03"
$ws.Range("E46").Value = "  +3.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.18"
$ws.Range("E47").Value = "  +2.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0532"
$ws.Range("E49").Value = "  +7.28%  "

$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.39"
$ws.Range("E51").Value = "  -0.44%  "
